# WikidataEnricher: resolve related items to get a full type hierarchy
#
# This script applies targeted cell-value updates to the synthetic_data
# workbook, touching the CmsWork, CmsWorkClosing, CmsWorkOpening, and
# CmsRightsStatement sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# CmsWork sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CmsWork")

$ws.Range("C2").Value = "_:Nc567c036f99a41cc9e0148c633049a60"
$ws.Range("D2").Value = "http://example.com/organization4"
$ws.Range("E2").Value = "CmsCollection0CmsWork1 alternative title 0"
$ws.Range("G2").Value = "CmsCollection0CmsWork1Id0"
$ws.Range("J2").Value = "CmsCollection0CmsWork1 provenance 1"
$ws.Range("U2").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:102"

$ws.Range("C3").Value = "_:Nedcadcc2e793425588240c8a1bcceb55"
$ws.Range("D3").Value = "http://example.com/person0"
$ws.Range("G3").Value = "CmsCollection0CmsWork3Id1"

$ws.Range("C4").Value = "_:N73d1f755c85c4a928483f10e2b054ca6"
$ws.Range("E4").Value = "CmsCollection1CmsWork5 alternative title 1"
$ws.Range("U4").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:105"

$ws.Range("C5").Value = "_:N032db3ea94f74aa2b4ee57d73173304b"
$ws.Range("G5").Value = "CmsCollection1CmsWork7Id0"
$ws.Range("J5").Value = "CmsCollection1CmsWork7 provenance 0"

$ws.Range("B6").Value = "_:Nfbe4e899e61c4e579f3f7e9a7b1ebb83"
$ws.Range("C6").Value = "http://example.com/organization2"
$ws.Range("D6").Value = "FreestandingWork9 alternative title 0"
$ws.Range("F6").Value = "FreestandingWork9Id1"
$ws.Range("I6").Value = "FreestandingWork9 provenance 1"
$ws.Range("T6").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:109"

$ws.Range("B7").Value = "_:N0a8e5eb4dc364da599fbcb23894ec198"
$ws.Range("C7").Value = "http://example.com/organization4"
$ws.Range("F7").Value = "FreestandingWork11Id1"

# ---------------------------------------------------------------------
# CmsWorkClosing sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CmsWorkClosing")

$ws.Range("A2").Value = "_:N4254e29086664f0dae1d3e48ec24d1b6"
$ws.Range("C2").Value = "_:Nddd06bc9340e4778a457a042dec15d34"

$ws.Range("A3").Value = "_:N33287ab2def7420a87360d4312c5bdf8"
$ws.Range("C3").Value = "_:Nc64014a4cde949deaa53d01f419cdf55"

$ws.Range("A4").Value = "_:Na04e302138b04f3dade7d80b90f615ce"
$ws.Range("C4").Value = "_:N25d5449e22ac488da37cc39f04594b76"

$ws.Range("A5").Value = "_:Ncf2e862331b14c5d8b93e913f91c616e"
$ws.Range("C5").Value = "_:N74d0b54ca09441879e888edc45406982"

$ws.Range("A6").Value = "_:N9afd0e9d078d404d88d58fceae853d47"
$ws.Range("C6").Value = "_:N12440732f0654826957808bebd87879b"

$ws.Range("A7").Value = "_:N1cce929b109646fb8d8e95680e8b6f1d"
$ws.Range("C7").Value = "_:N4ff4b16f678e4bb19d75d236016c2782"

# ---------------------------------------------------------------------
# CmsWorkOpening sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CmsWorkOpening")

$ws.Range("C2").Value = "_:Nddd06bc9340e4778a457a042dec15d34"
$ws.Range("C3").Value = "_:Nc64014a4cde949deaa53d01f419cdf55"
$ws.Range("C4").Value = "_:N25d5449e22ac488da37cc39f04594b76"
$ws.Range("C5").Value = "_:N74d0b54ca09441879e888edc45406982"
$ws.Range("C6").Value = "_:N12440732f0654826957808bebd87879b"
$ws.Range("C7").Value = "_:N4ff4b16f678e4bb19d75d236016c2782"

# ---------------------------------------------------------------------
# CmsRightsStatement sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CmsRightsStatement")

$ws.Range("E2").Value = "You may find additional information about the copyright status of the Item on the website of the organization that has made the Item available."
